$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 183
$ws.Range("F5").Value = 992
$ws.Range("F7").Value = 2460
$ws.Range("F8").Value = 31
$ws.Range("F9").Value = 1201
$ws.Range("F12").Value = 883
$ws.Range("F13").Value = 1091
$ws.Range("F15").Value = 288
$ws.Range("F17").Value = 714
$ws.Range("F18").Value = 749
$ws.Range("F19").Value = 178
$ws.Range("F20").Value = 469
$ws.Range("F21").Value = 1098
$ws.Range("F23").Value = 559
$ws.Range("F24").Value = 572
$ws.Range("F25").Value = 212
$ws.Range("F26").Value = 290
$ws.Range("F27").Value = 287
$ws.Range("F28").Value = 669
$ws.Range("F29").Value = 3089
$ws.Range("F33").Value = 26
$ws.Range("F35").Value = 124
$ws.Range("F36").Value = 1573
$ws.Range("F37").Value = 434
$ws.Range("F40").Value = 132
$ws.Range("F43").Value = 117
$ws.Range("F44").Value = 122
$ws.Range("F45").Value = 81

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 98
$ws.Range("F10").Value = 175
$ws.Range("F11").Value = 4411
$ws.Range("F13").Value = 23
$ws.Range("F14").Value = 167

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 710
$ws.Range("F4").Value = 664

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 710
$ws.Range("F5").Value = 664
$ws.Range("F7").Value = 183
$ws.Range("F8").Value = 992
$ws.Range("F9").Value = 2460
$ws.Range("F10").Value = 31
$ws.Range("F11").Value = 1201
$ws.Range("F14").Value = 883
$ws.Range("F15").Value = 1091
$ws.Range("F16").Value = 288
$ws.Range("F19").Value = 714
$ws.Range("F22").Value = 749
$ws.Range("F23").Value = 178
$ws.Range("F24").Value = 469
$ws.Range("F25").Value = 1098
$ws.Range("F26").Value = 98
$ws.Range("F28").Value = 559
$ws.Range("F29").Value = 572
$ws.Range("F30").Value = 212
$ws.Range("F31").Value = 287
$ws.Range("F33").Value = 3089
$ws.Range("F34").Value = 175
$ws.Range("F37").Value = 26
$ws.Range("F38").Value = 124
$ws.Range("F39").Value = 1573
$ws.Range("F40").Value = 434
$ws.Range("F44").Value = 132
$ws.Range("F46").Value = 117
$ws.Range("F47").Value = 81
